$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank columns at K:L. This shifts the existing Lat/Long
# (DMS-text) columns, and everything to their right, two columns to the right.
$ws.Range("K1:L1").EntireColumn.Insert()

# Header labels for the new decimal lat/lon columns.
$ws.Range("K1").Value = "Lat_dec"
$ws.Range("L1").Value = "Lon_dec"

# Decimal latitude/longitude values for each data row, converted from the
# degree-minute text already present (now in columns M/N).
$ws.Range("K2").Value = 26.5
$ws.Range("L2").Value = 96.58333

$ws.Range("K3").Value = 45.933
$ws.Range("L3").Value = -0.664

$ws.Range("K4").Value = 26.5
$ws.Range("L4").Value = 96.58333

$ws.Range("K5").Value = 26.5
$ws.Range("L5").Value = 96.58333

$ws.Range("K6").Value = 45.933
$ws.Range("L6").Value = -0.664

$ws.Range("K7").Value = 40.45
$ws.Range("L7").Value = -74.35

$ws.Range("K8").Value = 49.8166
$ws.Range("L8").Value = -111.6833

$ws.Range("K9").Value = 33.1333
$ws.Range("L9").Value = -35.58333

$ws.Range("K10").Value = 46.2333
$ws.Range("L10").Value = -2.68333

$ws.Range("K11").Value = 26.5
$ws.Range("L11").Value = 96.58333

# Resize the two new columns to fit their (short, numeric) content, like
# Excel would do on entry/auto-fit.
$ws.Columns.Item(11).ColumnWidth = 7.736979166666667
$ws.Columns.Item(12).ColumnWidth = 8.307291666666666

# Update view state: scroll so column J is the leftmost visible column and
# the last-entered cell (K11) is selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 10
$ws.Range("K11").Select()
